$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.430.64"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "1.643.79"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'212.04"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("E6").Value = "  +4.01%  "
$ws.Range("D8").Value = "'23.16"
$ws.Range("E8").Value = "  -3.07%  "
$ws.Range("E9").Value = "  -2.61%  "
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("D11").Value = "'0.0894"
$ws.Range("E11").Value = "  +1.75%  "
$ws.Range("D12").Value = "1.875.00"
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").Value = "1.640.88"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "'0.572"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("D15").Value = "'4.05"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").Value = "'64.43"
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("D17").Value = "27.404.43"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "'229.22"
$ws.Range("E18").Value = "  -5.17%  "
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").Value = "'7.59"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "'4.34"
$ws.Range("E22").Value = "  -3.32%  "
$ws.Range("D23").Value = "'9.60"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").Value = "'147.55"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "'15.61"
$ws.Range("E29").Value = "  -4.88%  "
$ws.Range("E30").Value = "  -3.57%  "
$ws.Range("E31").Value = "  -3.62%  "
$ws.Range("D32").Value = "'3.28"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("D33").Value = "'3.17"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("D34").Value = "1.418.00"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("E38").Value = "  -4.21%  "
$ws.Range("E39").Value = "  -4.13%  "
$ws.Range("D40").Value = "'1.02"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "'0.832"
$ws.Range("E42").Value = "  +5.42%  "
$ws.Range("D43").Value = "'2.47"
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").Value = "'64.59"
$ws.Range("E46").Value = "  -7.09%  "
$ws.Range("D47").Value = "1.784.70"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("E48").Value = "  -3.92%  "
$ws.Range("D49").Value = "'88.05"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("E51").Value = "  -3.05%  "
